$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Bills (row 3) CRUD operation flags
$ws.Range("B3").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1

# Update Receipts (row 4) CRUD operation flags
$ws.Range("B4").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1

# Update the active selection to match the new cursor position
$ws.Range("B4").Select()
